# Inserts a new price-report row for Achicoria (Vega Modelo de Temuco) above
# the existing row 44, shifting rows 44:85 down to 45:86, and fills the new
# row 44 with the week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 44; existing rows 44:85 shift to 45:86.
$ws.Rows.Item(44).Insert()

$row = 44
$ws.Cells.Item($row, 1).Value  = 10
$ws.Cells.Item($row, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value  = "La Araucanía"
$ws.Cells.Item($row, 4).Value  = 45040
$ws.Cells.Item($row, 5).Value  = 9
$ws.Cells.Item($row, 6).Value  = 100112010
$ws.Cells.Item($row, 7).Value  = "Achicoria"
$ws.Cells.Item($row, 8).Value  = "Sin especificar"
$ws.Cells.Item($row, 9).Value  = "Primera"
$ws.Cells.Item($row, 10).Value = 100
$ws.Cells.Item($row, 11).Value = 10000
$ws.Cells.Item($row, 12).Value = 10000
$ws.Cells.Item($row, 13).Value = 10000
$ws.Cells.Item($row, 14).Value = "`$/caja 18 unidades"
$ws.Cells.Item($row, 15).Value = "Región Metropolitana"
$ws.Cells.Item($row, 16).Value = 556
$ws.Cells.Item($row, 17).Value = 18
$ws.Cells.Item($row, 18).Value = "Hortaliza"
